$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim trailing spaces from header labels (row 1) in-place.
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "skills"
$ws.Range("E1").Value = "image"
$ws.Range("I1").Value = "date"

# Update the active selection to reflect the edit location (J1) as in the diff.
$ws.Range("J1").Select() | Out-Null
